# VUMI Global Flex benefits.xlsx - insert a space after the "$" placeholder
# in the various "Covered in full with $..." benefit descriptions, and
# update the saved selection to C17 (in progress - vumi flex).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# --- Diagnostics & Test / Surgeries & Anesthesia / Oncology (row 11-13) ---
$ws.Range("C11:G13").Value = "Covered in full with `$ deductible"

# --- Organ Transplant (row 14) ---
$ws.Range("C14:G14").Value = "Covered in full including USD 50,000 for donor cost with `$ deductible"

# --- Out-patient Consultations / Specialists (rows 16-17) ---
$ws.Range("C16:F17").Value = "Covered in full with `$ " + $nl

# --- Out-patient Medicines (row 18, columns C:E only - F18 is different text) ---
$ws.Range("C18:E18").Value = "Covered in full with `$ "

# --- Vaccination (row 19) - each column has its own USD amount ---
$ws.Range("C19").Value = "Travel and preventive vaccinations covered up to USD 500 with `$ "
$ws.Range("D19").Value = "Travel and preventive vaccinations covered up to USD 350 with `$ "
$ws.Range("E19").Value = "Travel and preventive vaccinations covered up to USD 250 with `$ "
$ws.Range("F19").Value = "Travel and preventive vaccinations covered up to USD 150 with `$ "

# --- Scans & Diagnostic Tests (row 20) ---
$ws.Range("C20:F20").Value = "Covered in full with `$ "

# --- Physiotherapy (row 21) ---
$ws.Range("C21:F21").Value = "Covered in full with `$ subject to pre-approval after 10 sessions" + $nl + "(Combined with Alternative medicines)"

# --- Alternative Medicines (row 33) ---
$ws.Range("C33:F33").Value = "Covered in full with `$ subject to pre-approval after 10 sessions" + $nl + "(Combined with Physiotherapy)"

# --- Mental Health Benefit (row 34, column C only) ---
$ws.Range("C34").Value = "OP Psychiatirc treatment covered up to USD 10,000 with `$ "

# --- Update the active selection left in the sheet view ---
$ws.Range("C17").Select()
